$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4d"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 68.50681433333334
$ws.Range("H2").Value = 205.520443
$ws.Range("I2").Value = 0.9663865053086182
$ws.Range("J2").Value = 0.9663865053086185
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.534538333333333
$ws.Range("N2").Value = 4.603615
$ws.Range("O2").Value = 0.1494637976135089
$ws.Range("P2").Value = 0.1494637976135089
$ws.Range("Q2").Value = 105.1263326890494
$ws.Range("R2").Value = 946.136994201445
$ws.Range("S2").Value = 0.1444397970458735
$ws.Range("T2").Value = 0.1444397970458735

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4d"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 68.50681433333334
$ws.Range("H3").Value = 205.520443
$ws.Range("I3").Value = 0.9663865053086182
$ws.Range("J3").Value = 0.9663865053086185
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.973328333333334
$ws.Range("N3").Value = 14.919985
$ws.Range("O3").Value = 0.4844014146353658
$ws.Range("P3").Value = 0.4844014146353658
$ws.Range("Q3").Value = 340.7068807503728
$ws.Range("R3").Value = 3066.361926753355
$ws.Range("S3").Value = 0.4681189902560222
$ws.Range("T3").Value = 0.4681189902560222

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4d"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 68.50681433333334
$ws.Range("H4").Value = 205.520443
$ws.Range("I4").Value = 0.9663865053086182
$ws.Range("J4").Value = 0.9663865053086185
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.75909
$ws.Range("N4").Value = 11.27727
$ws.Range("O4").Value = 0.3661347877511252
$ws.Range("P4").Value = 0.3661347877511252
$ws.Range("Q4").Value = 257.52328069229
$ws.Range("R4").Value = 2317.70952623061
$ws.Range("S4").Value = 0.3538277180067225
$ws.Range("T4").Value = 0.3538277180067226

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema4d"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.228643
$ws.Range("H5").Value = 3.685929
$ws.Range("I5").Value = 0.01733176511849816
$ws.Range("J5").Value = 0.01733176511849817
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 1.534538333333333
$ws.Range("N5").Value = 4.603615
$ws.Range("O5").Value = 0.1494637976135089
$ws.Range("P5").Value = 0.1494637976135089
$ws.Range("Q5").Value = 1.885399781481666
$ws.Range("R5").Value = 16.968598033335
$ws.Range("S5").Value = 0.002590471433956083
$ws.Range("T5").Value = 0.002590471433956084

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4d"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.228643
$ws.Range("H6").Value = 3.685929
$ws.Range("I6").Value = 0.01733176511849816
$ws.Range("J6").Value = 0.01733176511849817
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.973328333333334
$ws.Range("N6").Value = 14.919985
$ws.Range("O6").Value = 0.4844014146353658
$ws.Range("P6").Value = 0.4844014146353658
$ws.Range("Q6").Value = 6.110445043451667
$ws.Range("R6").Value = 54.994005391065
$ws.Range("S6").Value = 0.008395531541528398
$ws.Range("T6").Value = 0.0083955315415284

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4d"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.228643
$ws.Range("H7").Value = 3.685929
$ws.Range("I7").Value = 0.01733176511849816
$ws.Range("J7").Value = 0.01733176511849817
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.75909
$ws.Range("N7").Value = 11.27727
$ws.Range("O7").Value = 0.3661347877511252
$ws.Range("P7").Value = 0.3661347877511252
$ws.Range("Q7").Value = 4.61857961487
$ws.Range("R7").Value = 41.56721653383
$ws.Range("S7").Value = 0.006345762143013679
$ws.Range("T7").Value = 0.00634576214301368

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema4d"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.154206333333333
$ws.Range("H8").Value = 3.462619
$ws.Range("I8").Value = 0.01628172957288352
$ws.Range("J8").Value = 0.01628172957288353
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.534538333333333
$ws.Range("N8").Value = 4.603615
$ws.Range("O8").Value = 0.1494637976135089
$ws.Range("P8").Value = 0.1494637976135089
$ws.Range("Q8").Value = 1.771173863076111
$ws.Range("R8").Value = 15.940564767685
$ws.Range("S8").Value = 0.002433529133679346
$ws.Range("T8").Value = 0.002433529133679347

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema4d"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.154206333333333
$ws.Range("H9").Value = 3.462619
$ws.Range("I9").Value = 0.01628172957288352
$ws.Range("J9").Value = 0.01628172957288353
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.973328333333334
$ws.Range("N9").Value = 14.919985
$ws.Range("O9").Value = 0.4844014146353658
$ws.Range("P9").Value = 0.4844014146353658
$ws.Range("Q9").Value = 5.740247060079446
$ws.Range("R9").Value = 51.662223540715
$ws.Range("S9").Value = 0.007886892837815249
$ws.Range("T9").Value = 0.007886892837815251

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema4d"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.154206333333333
$ws.Range("H10").Value = 3.462619
$ws.Range("I10").Value = 0.01628172957288352
$ws.Range("J10").Value = 0.01628172957288353
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.75909
$ws.Range("N10").Value = 11.27727
$ws.Range("O10").Value = 0.3661347877511252
$ws.Range("P10").Value = 0.3661347877511252
$ws.Range("Q10").Value = 4.338765485570001
$ws.Range("R10").Value = 39.04888937013
$ws.Range("S10").Value = 0.005961307601388927
$ws.Range("T10").Value = 0.005961307601388928

Write-Output "done"